$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UE")

# Insert a new row before row 8, shifting the existing rows (old row 8 onward) down by one.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new SNMP commit command entry.
$ws.Cells.Item(8,1).Value = "pmpDevCpeLteBandCommit "
$ws.Cells.Item(8,2).Value = ".1.3.6.1.4.1.17713.20.2.3.1.2.15.0"
$ws.Cells.Item(8,3).Value = "commit"
$ws.Cells.Item(8,4).Value = "Oct"

# Match the style of the other "Zemtel" rows (style index 9, same font/format as 6/7/8).
$ws.Range("A8:D8").Style = $ws.Range("A14:D14").Style

$ws.Range("D9").Select()
